# Add color palette slide — reposition the two existing swatch rectangles on
# slide 4 and add the remaining swatches ("Sliders", stack/pie trio, line
# trio, and the two left-column cards) described in the commit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# --- helpers -----------------------------------------------------------
# EMU -> points conversion for Shapes.AddShape (which already round-trips
# exactly), and for the Left/Top setters (which truncate to 4dp before
# converting back to EMU, so a tiny epsilon keeps them on target).
function EmuToPt($emu) {
    return $emu / 12700.0
}
function EmuToPtForSetter($emu) {
    return ($emu / 12700.0) + 0.00005
}
# VBA-style RGB(): low byte = R, mid byte = G, high byte = B.
function RGBColor($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}
function HexToRGBColor($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return RGBColor $r $g $b
}

# --- move the two existing swatches -------------------------------------
$bg = $s.Shapes.Item("Rectangle 5")
$bg.Left = EmuToPtForSetter 1465219
$bg.Top = EmuToPtForSetter 1792343

$chartsBg = $s.Shapes.Item("Rectangle 6")
$chartsBg.Left = EmuToPtForSetter 1465219
$chartsBg.Top = EmuToPtForSetter 2692675

# --- add the new swatches ------------------------------------------------
function New-Swatch($name, $x, $y, $cx, $cy, $fillHex, $noLine) {
    $shp = $s.Shapes.AddShape(1, (EmuToPt $x), (EmuToPt $y), (EmuToPt $cx), (EmuToPt $cy))
    $shp.Name = $name
    $shp.Fill.ForeColor.RGB = HexToRGBColor $fillHex
    if ($noLine) {
        $shp.Line.Visible = $false
    }
    $tf = $shp.TextFrame
    $tf.VerticalAnchor = 3
    $tr = $tf.TextRange
    $tr.ParagraphFormat.Alignment = 2
    return $shp
}

$sliders = New-Swatch "Rectangle 7" 7196317 1792343 3812344 717452 "1D3D70" $false
$sliders.TextFrame.TextRange.Text = "Sliders"

$bottomStack = New-Swatch "Rectangle 8" 7196317 2692675 3812344 717452 "054380" $false
$bottomStack.TextFrame.TextRange.Text = " "
$bottomStack.TextFrame.TextRange.InsertAfter("Bottom Stack / Pie") | Out-Null

$middleStack = New-Swatch "Rectangle 9" 7196317 3593007 3812344 717452 "317ABE" $false
$middleStack.TextFrame.TextRange.Text = " "
$middleStack.TextFrame.TextRange.InsertAfter("Middle Stack / Pie") | Out-Null

$topStack = New-Swatch "Rectangle 10" 7196317 4493339 3812344 717452 "6FB3F2" $true
$topStack.TextFrame.TextRange.Text = " "
$topStack.TextFrame.TextRange.InsertAfter("Top Stack / Pie") | Out-Null

$dayCase = New-Swatch "Rectangle 11" 7196317 5393671 1215998 717452 "E96D1A" $true
$dayCase.TextFrame.TextRange.Text = "Day Case Line"

$inpatient = New-Swatch "Rectangle 12" 8494490 5393671 1215998 717452 "D94451" $true
$inpatient.TextFrame.TextRange.Text = "Inpatient Line"

$outpatient = New-Swatch "Rectangle 13" 9792663 5393671 1215998 717452 "078CFF" $true
$outpatient.TextFrame.TextRange.Text = "Outpatient Line"

$cardHeadline = New-Swatch "Rectangle 14" 1465219 3593007 3812344 717452 "8DC0EC" $true
$cardHeadline.TextFrame.TextRange.Text = "Card Headline"

$boxOutline = New-Swatch "Rectangle 15" 1465219 4493339 3812344 717452 "3A57B4" $true
$boxOutline.TextFrame.TextRange.Text = "Box Outline"
